$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.178.63'
$ws.Range("E2").Value = '  +3.91%  '

$ws.Range("D3").Value = '1.780.47'
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.11'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9989'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3822'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3437'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.61%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.47'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.96%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.156'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07422'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.30'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +7.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9970'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.432'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.02%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.788.26'
$ws.Range("E15").Value = '  +0.33%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.177'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001080'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06661'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.43'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9986'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.49'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.450'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.90%  '

$ws.Range("D23").Value = '28.286.54'
$ws.Range("E23").Value = '  +4.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.13'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.373'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.89'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.437'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.416'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '155.22'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("E30").Value = '  +1.64%  '

$ws.Range("D31").Value = '1.988.98'
$ws.Range("E31").Value = '  +0.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.150'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.964'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08885'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.80'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.97%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02437'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6881'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.334'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06351'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2184'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.241'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.43%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.496'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.324'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.14%  '

$ws.Range("E44").Value = '  +0.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9979'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.26%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6303'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.867'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.95'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.101'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07458'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.214'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +9.23%  '
